# Remove the last three slides ("Première normalisation", "Segmentation",
# "Deuxième normalisation") from the deck, per the commit:
# only slides 1-3 (OncoscanR title, Input, Pre-process) remain.

$p = $ppt.ActivePresentation

$p.Slides.Item(6).Delete()
$p.Slides.Item(5).Delete()
$p.Slides.Item(4).Delete()

Write-Output $p.Slides.Count
